$wb = $excel.ActiveWorkbook

# Switch to the PVTStL sheet (it becomes the active/selected tab)
$ws = $wb.Worksheets.Item("PVTStL")
$ws.Activate()

# "ships" (row 6) is now marked as subject to LCFS for both passenger and freight
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# Leave the selection on B6, matching the saved view state
$ws.Range("B6").Select()
